# Updated data and added a new "Loss of Smell or Taste" symptom row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 ("Shortness of Breath"),
# shifting the old rows 4 and 5 down to 5 and 6.
$ws.Rows.Item(4).Insert()

# Row 3 (Position 2) now represents the new "Loss of Smell or Taste" symptom.
$ws.Range("D3").Value = "Loss of Smell or Taste"

# The newly inserted row 4 (Position 3) holds the updated "Cough" data.
$ws.Range("B4").Value = 4.5
$ws.Range("C4").Value = 7.5
$ws.Range("D4").Value = "Cough"

# Renumber the Position column (A) so it stays sequential 1..5.
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Match the saved selection state.
$ws.Range("E4").Select()
